$wb = $excel.ActiveWorkbook

$wsExclude = $wb.Worksheets.Item("Exclude")
$wsInclude = $wb.Worksheets.Item("Include")

# ------------------------------------------------------------------
# Rename the "brand" parameter to "brand_name" on the Include sheet
# (Numerator / Denominator rows that key off the brand param), part of
# widening the white-milk brand exclusion.
# ------------------------------------------------------------------
$wsInclude.Range("C2").Value = "brand_name"
$wsInclude.Range("C4").Value = "brand_name"

# ------------------------------------------------------------------
# Normalise the font on the "White milk" exclusion rows (D3, D5, D6,
# D8) of the Exclude sheet so they match the formatting already used
# by the other Value-1 cells (D2, D4, D7, D9) -- this collapses the
# redundant duplicate font/style that only those four cells used.
# ------------------------------------------------------------------
foreach ($addr in @("D3", "D5", "D6", "D8")) {
    $cell = $wsExclude.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Color = 0
}

# ------------------------------------------------------------------
# Switch the active/selected sheet from Exclude to Include, and move
# the view/selection on Include down to the newly relevant rows.
# ------------------------------------------------------------------
$wsInclude.Activate()
$wsInclude.Range("E30").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
